$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update time_taken column (F2:F131) with refreshed timestamps
$ws.Range("F2").Value = "2021-10-05 14:33:40.010156"
$ws.Range("F3").Value = "2021-10-05 14:33:40.010164"
$ws.Range("F4").Value = "2021-10-05 14:33:40.010167"
$ws.Range("F5").Value = "2021-10-05 14:33:40.010170"
$ws.Range("F6").Value = "2021-10-05 14:33:40.010173"
$ws.Range("F7").Value = "2021-10-05 14:33:40.010175"
$ws.Range("F8").Value = "2021-10-05 14:33:40.010178"
$ws.Range("F9").Value = "2021-10-05 14:33:40.010180"
$ws.Range("F10").Value = "2021-10-05 14:33:40.010183"
$ws.Range("F11").Value = "2021-10-05 14:33:40.010185"
$ws.Range("F12").Value = "2021-10-05 14:33:40.010188"
$ws.Range("F13").Value = "2021-10-05 14:33:40.010191"
$ws.Range("F14").Value = "2021-10-05 14:33:40.010193"
$ws.Range("F15").Value = "2021-10-05 14:33:40.010196"
$ws.Range("F16").Value = "2021-10-05 14:33:40.010198"
$ws.Range("F17").Value = "2021-10-05 14:33:40.010200"
$ws.Range("F18").Value = "2021-10-05 14:33:40.010203"
$ws.Range("F19").Value = "2021-10-05 14:33:40.010206"
$ws.Range("F20").Value = "2021-10-05 14:33:40.010208"
$ws.Range("F21").Value = "2021-10-05 14:33:40.010211"
$ws.Range("F22").Value = "2021-10-05 14:33:40.010213"
$ws.Range("F23").Value = "2021-10-05 14:33:40.010216"
$ws.Range("F24").Value = "2021-10-05 14:33:40.010218"
$ws.Range("F25").Value = "2021-10-05 14:33:40.010221"
$ws.Range("F26").Value = "2021-10-05 14:33:40.010224"
$ws.Range("F27").Value = "2021-10-05 14:33:40.010226"
$ws.Range("F28").Value = "2021-10-05 14:33:40.010229"
$ws.Range("F29").Value = "2021-10-05 14:33:40.010232"
$ws.Range("F30").Value = "2021-10-05 14:33:40.010234"
$ws.Range("F31").Value = "2021-10-05 14:33:40.010237"
$ws.Range("F32").Value = "2021-10-05 14:33:40.010239"
$ws.Range("F33").Value = "2021-10-05 14:33:40.010242"
$ws.Range("F34").Value = "2021-10-05 14:33:40.010245"
$ws.Range("F35").Value = "2021-10-05 14:33:40.010247"
$ws.Range("F36").Value = "2021-10-05 14:33:40.010249"
$ws.Range("F37").Value = "2021-10-05 14:33:40.010252"
$ws.Range("F38").Value = "2021-10-05 14:33:40.010254"
$ws.Range("F39").Value = "2021-10-05 14:33:40.010257"
$ws.Range("F40").Value = "2021-10-05 14:33:40.010259"
$ws.Range("F41").Value = "2021-10-05 14:33:40.010262"
$ws.Range("F42").Value = "2021-10-05 14:33:40.010265"
$ws.Range("F43").Value = "2021-10-05 14:33:40.010267"
$ws.Range("F44").Value = "2021-10-05 14:33:40.010270"
$ws.Range("F45").Value = "2021-10-05 14:33:40.010272"
$ws.Range("F46").Value = "2021-10-05 14:33:40.010275"
$ws.Range("F47").Value = "2021-10-05 14:33:40.010277"
$ws.Range("F48").Value = "2021-10-05 14:33:40.010280"
$ws.Range("F49").Value = "2021-10-05 14:33:40.010282"
$ws.Range("F50").Value = "2021-10-05 14:33:40.010284"
$ws.Range("F51").Value = "2021-10-05 14:33:40.010287"
$ws.Range("F52").Value = "2021-10-05 14:33:40.010289"
$ws.Range("F53").Value = "2021-10-05 14:33:40.010292"
$ws.Range("F54").Value = "2021-10-05 14:33:40.010295"
$ws.Range("F55").Value = "2021-10-05 14:33:40.010297"
$ws.Range("F56").Value = "2021-10-05 14:33:40.010300"
$ws.Range("F57").Value = "2021-10-05 14:33:40.010302"
$ws.Range("F58").Value = "2021-10-05 14:33:40.010305"
$ws.Range("F59").Value = "2021-10-05 14:33:40.010307"
$ws.Range("F60").Value = "2021-10-05 14:33:40.010310"
$ws.Range("F61").Value = "2021-10-05 14:33:40.010312"
$ws.Range("F62").Value = "2021-10-05 14:33:40.010314"
$ws.Range("F63").Value = "2021-10-05 14:33:40.010317"
$ws.Range("F64").Value = "2021-10-05 14:33:40.010319"
$ws.Range("F65").Value = "2021-10-05 14:33:40.010322"
$ws.Range("F66").Value = "2021-10-05 14:33:40.010325"
$ws.Range("F67").Value = "2021-10-05 14:33:40.010328"
$ws.Range("F68").Value = "2021-10-05 14:33:40.010331"
$ws.Range("F69").Value = "2021-10-05 14:33:40.010333"
$ws.Range("F70").Value = "2021-10-05 14:33:40.010336"
$ws.Range("F71").Value = "2021-10-05 14:33:40.010338"
$ws.Range("F72").Value = "2021-10-05 14:33:40.010341"
$ws.Range("F73").Value = "2021-10-05 14:33:40.010343"
$ws.Range("F74").Value = "2021-10-05 14:33:40.010346"
$ws.Range("F75").Value = "2021-10-05 14:33:40.010348"
$ws.Range("F76").Value = "2021-10-05 14:33:40.010351"
$ws.Range("F77").Value = "2021-10-05 14:33:40.010353"
$ws.Range("F78").Value = "2021-10-05 14:33:40.010358"
$ws.Range("F79").Value = "2021-10-05 14:33:40.010361"
$ws.Range("F80").Value = "2021-10-05 14:33:40.010364"
$ws.Range("F81").Value = "2021-10-05 14:33:40.010366"
$ws.Range("F82").Value = "2021-10-05 14:33:40.010369"
$ws.Range("F83").Value = "2021-10-05 14:33:40.010371"
$ws.Range("F84").Value = "2021-10-05 14:33:40.010374"
$ws.Range("F85").Value = "2021-10-05 14:33:40.010376"
$ws.Range("F86").Value = "2021-10-05 14:33:40.010379"
$ws.Range("F87").Value = "2021-10-05 14:33:40.010381"
$ws.Range("F88").Value = "2021-10-05 14:33:40.010384"
$ws.Range("F89").Value = "2021-10-05 14:33:40.010386"
$ws.Range("F90").Value = "2021-10-05 14:33:40.010389"
$ws.Range("F91").Value = "2021-10-05 14:33:40.010391"
$ws.Range("F92").Value = "2021-10-05 14:33:40.010394"
$ws.Range("F93").Value = "2021-10-05 14:33:40.010396"
$ws.Range("F94").Value = "2021-10-05 14:33:40.010400"
$ws.Range("F95").Value = "2021-10-05 14:33:40.010403"
$ws.Range("F96").Value = "2021-10-05 14:33:40.010406"
$ws.Range("F97").Value = "2021-10-05 14:33:40.010408"
$ws.Range("F98").Value = "2021-10-05 14:33:40.010411"
$ws.Range("F99").Value = "2021-10-05 14:33:40.010413"
$ws.Range("F100").Value = "2021-10-05 14:33:40.010416"
$ws.Range("F101").Value = "2021-10-05 14:33:40.010418"
$ws.Range("F102").Value = "2021-10-05 14:33:40.010421"
$ws.Range("F103").Value = "2021-10-05 14:33:40.010423"
$ws.Range("F104").Value = "2021-10-05 14:33:40.010426"
$ws.Range("F105").Value = "2021-10-05 14:33:40.010428"
$ws.Range("F106").Value = "2021-10-05 14:33:40.010431"
$ws.Range("F107").Value = "2021-10-05 14:33:40.010433"
$ws.Range("F108").Value = "2021-10-05 14:33:40.010436"
$ws.Range("F109").Value = "2021-10-05 14:33:40.010438"
$ws.Range("F110").Value = "2021-10-05 14:33:40.010443"
$ws.Range("F111").Value = "2021-10-05 14:33:40.010446"
$ws.Range("F112").Value = "2021-10-05 14:33:40.010448"
$ws.Range("F113").Value = "2021-10-05 14:33:40.010451"
$ws.Range("F114").Value = "2021-10-05 14:33:40.010453"
$ws.Range("F115").Value = "2021-10-05 14:33:40.010456"
$ws.Range("F116").Value = "2021-10-05 14:33:40.010458"
$ws.Range("F117").Value = "2021-10-05 14:33:40.010461"
$ws.Range("F118").Value = "2021-10-05 14:33:40.010463"
$ws.Range("F119").Value = "2021-10-05 14:33:40.010466"
$ws.Range("F120").Value = "2021-10-05 14:33:40.010468"
$ws.Range("F121").Value = "2021-10-05 14:33:40.010470"
$ws.Range("F122").Value = "2021-10-05 14:33:40.010473"
$ws.Range("F123").Value = "2021-10-05 14:33:40.010475"
$ws.Range("F124").Value = "2021-10-05 14:33:40.010478"
$ws.Range("F125").Value = "2021-10-05 14:33:40.010480"
$ws.Range("F126").Value = "2021-10-05 14:33:40.010483"
$ws.Range("F127").Value = "2021-10-05 14:33:40.010485"
$ws.Range("F128").Value = "2021-10-05 14:33:40.010488"
$ws.Range("F129").Value = "2021-10-05 14:33:40.010490"
$ws.Range("F130").Value = "2021-10-05 14:33:40.010495"
$ws.Range("F131").Value = "2021-10-05 14:33:40.010498"

# Add the new "metadata" sheet positioned after "data"
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$metaSheet.Name = "metadata"

# Header row (bold, centered, bordered - copy formatting from the "data" header cells)
$ws.Range("B1").Copy() | Out-Null
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Deafness_Isolated"
$metaSheet.Range("C2").Value = 3241
# "1.13" must stay text (not be coerced to the number 1.13) and keep the
# default (unstyled) cell format, so force text via NumberFormat, assign,
# then clear the format back off again (keeps the stored type as text).
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.13"
$metaSheet.Range("D2").ClearFormats()
$metaSheet.Range("E2").Value = "2021-08-25T23:32:46.576427Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:40.006462"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3241/?format=json"

# Restore "data" as the active sheet/selection (only the sheet list changed)
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
